$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 and B3 values
$ws.Range("B2").Value = 2106
$ws.Range("B3").Value = 1608

# Row 4 now takes on the values previously held in row 5 (A=2) but with a new B value
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 544

# Remove the old row 5 entirely (shifts rows up, nothing below so it just clears/removes it)
$ws.Range("A5:B5").Delete()
